# Commit: "updates several machine setups"
#
# The "Commands" sheet had a duplicated/bogus row 20
# (both B20 and C20 held the same stray string "writes values to the
# registers in slaves specified by the given id" instead of real
# Command/Documentation text). That row is removed entirely, which shifts
# every row below it (21..91) up by one (20..90).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Delete the whole row 20 - everything below shifts up to fill the gap.
$ws.Rows(20).EntireRow.Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("A20").Select() | Out-Null
